$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.979.30'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '3.113.37'
$ws.Range("E3").Value = '  +2.59%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'579.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").Value = "'173.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.78%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.108.21'
$ws.Range("E8").Value = '  +2.59%  '
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("D10").Value = "'6.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.73%  '
$ws.Range("E11").Value = '  +1.61%  '
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("D14").Value = "'37.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '3.627.82'
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("D17").Value = '67.012.28'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").Value = '3.113.82'
$ws.Range("E19").Value = '  +2.55%  '
$ws.Range("D20").Value = "'16.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.08%  '
$ws.Range("D21").Value = "'487.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.83%  '
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("D23").Value = "'7.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.04%  '
$ws.Range("D24").Value = "'84.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.60%  '
$ws.Range("D25").Value = "'13.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.44%  '
$ws.Range("E26").Value = '  +4.32%  '
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("D29").Value = "'8.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.15%  '
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("E31").Value = '  +2.60%  '
$ws.Range("E32").Value = '  +3.22%  '
$ws.Range("D33").Value = "'0.0000100"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("E34").Value = '  -4.08%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").Value = "'47.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.52%  '
$ws.Range("E39").Value = '  +3.26%  '
$ws.Range("D40").Value = "'50.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.13%  '
$ws.Range("E41").Value = '  +2.46%  '
$ws.Range("E42").Value = '  +1.62%  '
$ws.Range("E43").Value = '  +0.74%  '
$ws.Range("E44").Value = '  -1.46%  '
$ws.Range("D45").Value = '2.848.95'
$ws.Range("E45").Value = '  +4.34%  '
$ws.Range("D46").Value = "'385.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").Value = "'136.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("D50").Value = "'25.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("E51").Value = '  +0.42%  '
